# Controle dos bancos de dados
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pesquisa Rolf")

# --- New header cell L4 ---
$ws.Range("L4").Value = "Extraído csv"

# --- New J column values (rows 6,9,10,11,14,15), matching style of existing J cells ---
$ws.Range("J6").Value = 11511559
$ws.Range("J9").Value = 11429889
$ws.Range("J10").Value = 11024228
$ws.Range("J11").Value = 11612715
$ws.Range("J14").Value = 11675269
$ws.Range("J15").Value = 8966476

foreach ($addr in @("J6","J9","J10","J11","J14","J15")) {
    $ws.Range($addr).NumberFormat = "#,##0"
}

# --- Column K formulas: shared formula over K5:K11, individual formulas K14 & K15 ---
$ws.Range("K5:K11").Formula = '=IF(J5=F5,"OK","DEU RUIM")'
$ws.Range("K14").Formula = '=IF(J14=F14,"OK","DEU RUIM")'
$ws.Range("K15").Formula = '=IF(J15=F15,"OK","DEU RUIM")'

# --- Column L "SIM" for rows whose OK result is confirmed ---
$ws.Range("L6").Value = "SIM"
$ws.Range("L7").Value = "SIM"
$ws.Range("L8").Value = "SIM"
$ws.Range("L11").Value = "SIM"
$ws.Range("L14").Value = "SIM"

# --- Conditional formatting on K5:K15: highlight cells containing "RUIM" ---
$rng = $ws.Range("K5:K15")
$fc = $rng.FormatConditions.Add(9, 0, "RUIM")
$fc.Text = "RUIM"
$fc.Formula1 = 'NOT(ISERROR(SEARCH("RUIM",K5)))'
$fc.Font.Bold = $true
$fc.Font.Color = 255

# --- Size new columns K (11) and L (12) to match new content width ---
$ws.Columns("K").ColumnWidth = 8.9
$ws.Columns("L").ColumnWidth = 10.45

# --- Selection ---
$ws.Range("L5").Select() | Out-Null
